# Replace "kubectl create -f <file>" with "kubectl apply -f <file>" on the
# three "kubectl create <object>.yaml" demo command lines (Service, ConfigMap,
# Secret slides). Also re-splits two sibling "kubectl create ..." command
# lines (configmap / secret generic) into the same finer-grained run layout,
# even though their text content is unchanged, to mirror how PowerPoint
# reflows runs when a shape's text is touched.

$p = $ppt.ActivePresentation
$enDash = [char]0x2013

function Edit-CreateToApply {
    param($para)
    # before: "$ kubectl create <EN-DASH>f my-XXXX.yaml"
    # after : "$ kubectl apply  <EN-DASH>f my-XXXX.yaml"
    $para.Characters(1,2).Text = "`$ "
    $para.Characters(3,7).Text = "kubectl"
    $para.Characters(10,1).Text = " "
    $para.Characters(11,7).Text = "apply "
    $para.Characters(17,3).Text = "$enDash" + "f "
}

# ---------------------------------------------------------------------
# Slide 16 ("Service"): shape "CustomShape 7", paragraph 2
#   "$ kubectl create -f my-service.yaml"  ->  "... apply ... my-service.yaml"
# ---------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(7)
$para16 = $sh16.TextFrame.TextRange.Paragraphs(2)
Edit-CreateToApply $para16
$para16.Characters(20,3).Text = "my-"
$para16.Characters(23,12).Text = "service.yaml"

# ---------------------------------------------------------------------
# Slide 28 ("ConfigMap"): shape "CustomShape 2"
#   Paragraph 1: "$ kubectl create -f my-config.yaml" -> apply
#   Paragraph 2: "$ kubectl create configmap special-config --from-literal=special.how=very"
#                (text unchanged, only re-split into more runs)
# ---------------------------------------------------------------------
$s28 = $p.Slides.Item(28)
$sh28 = $s28.Shapes.Item(4)
$tr28 = $sh28.TextFrame.TextRange

$para28_1 = $tr28.Paragraphs(1)
Edit-CreateToApply $para28_1
$para28_1.Characters(20,3).Text = "my-"
$para28_1.Characters(23,11).Text = "config.yaml"

$para28_2 = $tr28.Paragraphs(2)
$para28_2.Characters(1,2).Text = "`$ "
$para28_2.Characters(3,7).Text = "kubectl"
$para28_2.Characters(10,8).Text = " create "
$para28_2.Characters(18,9).Text = "configmap"
$para28_2.Characters(27,31).Text = " special-config --from-literal="
$para28_2.Characters(58,11).Text = "special.how"
$para28_2.Characters(69,5).Text = "=very"

# ---------------------------------------------------------------------
# Slide 30 ("Secret"): shape "CustomShape 3"
#   Paragraph 1: "$ kubectl create -f my-secret.yaml" -> apply
#   Paragraph 2: "$ kubectl create secret generic mysecret --from-literal=..."
#                (text unchanged, only re-split into more runs)
# ---------------------------------------------------------------------
$s30 = $p.Slides.Item(30)
$sh30 = $s30.Shapes.Item(4)
$tr30 = $sh30.TextFrame.TextRange

$para30_1 = $tr30.Paragraphs(1)
Edit-CreateToApply $para30_1
$para30_1.Characters(20,3).Text = "my-"
$para30_1.Characters(23,11).Text = "secret.yaml"

$para30_2 = $tr30.Paragraphs(2)
$para30_2.Characters(1,2).Text = "`$ "
$para30_2.Characters(3,7).Text = "kubectl"
$para30_2.Characters(10,23).Text = " create secret generic "
$para30_2.Characters(33,8).Text = "mysecret"
$para30_2.Characters(41,68).Text = " --from-literal=username=admin --from-literal=password=1f2d1e2e67df "
